# edit.ps1 -- applies the Concertacionevaluacion.docx content edits:
#   1. "Quices, ejercicios, entregas parciales trabajo final" -> the word
#      "Quices" is split into its own run and wrapped with
#      <w:proofErr w:type="spellStart"/> / <w:proofErr w:type="spellEnd"/>
#      (the rest of the sentence becomes a second, separately-formatted run).
#   2. Remove the yellow <w:highlight/> from the two "9.00 AM - 11.00 AM"
#      runs and from "MIERCOLES" / "VIERNES".
#   3. Change the day name "MIERCOLES" to "LUNES".
#
# Because these target runs share identical character formatting with
# their neighbours, the host's high level Range.InsertAfter / Range.Text
# writers silently coalesce adjacent runs, and Range.HighlightColorIndex
# does not persist through to the underlying <w:highlight/> element, so
# each edit below is performed by replacing the whole owning paragraph
# with freshly authored OOXML via Range.InsertXML - that is the one
# primitive that reliably round-trips exact run/proofErr structure.
#
# All five target paragraphs appear once each, in this order, top to
# bottom in the document, so a single forward-moving search cursor
# (rather than re-searching from the top every time, which would just
# keep re-matching the first "9.00 AM - 11.00 AM") visits each of them
# exactly once.

$d = $word.ActiveDocument
$cursor = 0

function Replace-NextParagraph($findText, $paragraphXml) {
    $script:cursor
    $r = $d.Range($cursor, $d.Content.End)
    $ok = $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Replace-NextParagraph: text not found after cursor $cursor : $findText"
    }
    $matchStart = $r.Start
    $r.Expand(4) | Out-Null
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/part.xml" pkg:contentType="application/xml"><pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:body>' + $paragraphXml + '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($pkg)
    $script:cursor = $matchStart + 1
}

# --- 1. Quices, ejercicios... -> split "Quices" into its own proofErr-wrapped run
$quicesPara = '<w:p w14:paraId="516CB9B2" w14:textId="09E4A85D" w:rsidR="0016663E" w:rsidRPr="00643EEB" w:rsidRDefault="0016663E" w:rsidP="00B65420">' +
    '<w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:lastRenderedPageBreak/><w:t>Quices</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>, ejercicios, entregas parciales trabajo final</w:t></w:r>' +
    '</w:p>'
Replace-NextParagraph "Quices, ejercicios, entregas parciales trabajo final" $quicesPara

# --- 2/3. Schedule table: drop the yellow highlight (two "9.00 AM - 11.00 AM"
#          runs, "MIERCOLES" and "VIERNES"), and rename MIERCOLES -> LUNES.
$firstTimePara = '<w:p w14:paraId="3086D333" w14:textId="497F548C" w:rsidR="003F0259" w:rsidRPr="00747BC9" w:rsidRDefault="00CF3CC2" w:rsidP="00B91A92">' +
    '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr>' +
    '<w:r w:rsidRPr="002634EF"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>9.00 AM - 11.00 AM</w:t></w:r>' +
    '</w:p>'
Replace-NextParagraph "9.00 AM - 11.00 AM" $firstTimePara

$lunesPara = '<w:p w14:paraId="24AAE5A5" w14:textId="30911388" w:rsidR="003F0259" w:rsidRPr="00747BC9" w:rsidRDefault="0085267B" w:rsidP="00B91A92">' +
    '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr>' +
    '<w:r w:rsidRPr="002634EF"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>LUNES</w:t></w:r>' +
    '</w:p>'
Replace-NextParagraph "MIERCOLES" $lunesPara

$secondTimePara = '<w:p w14:paraId="5B4E6CB0" w14:textId="56AA194E" w:rsidR="003F0259" w:rsidRPr="00747BC9" w:rsidRDefault="002634EF" w:rsidP="00B91A92">' +
    '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr>' +
    '<w:r w:rsidRPr="002634EF"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>9.00 AM - 11.00 AM</w:t></w:r>' +
    '</w:p>'
Replace-NextParagraph "9.00 AM - 11.00 AM" $secondTimePara

$viernesPara = '<w:p w14:paraId="3C5D860F" w14:textId="70FF7214" w:rsidR="003F0259" w:rsidRPr="00747BC9" w:rsidRDefault="002634EF" w:rsidP="00B91A92">' +
    '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>VIERNES</w:t></w:r>' +
    '</w:p>'
Replace-NextParagraph "VIERNES" $viernesPara

Write-Output "edits applied"
